$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.460980666666667
$ws.Range("H2").Value = 7.382942000000001
$ws.Range("I2").Value = 0.3244396275423151
$ws.Range("J2").Value = 0.3244396275423151
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4445023333333333
$ws.Range("N2").Value = 1.333507
$ws.Range("O2").Value = 0.009977046255258984
$ws.Range("P2").Value = 0.009977046255258982
$ws.Range("Q2").Value = 1.093911648621556
$ws.Range("R2").Value = 9.845204837594
$ws.Range("S2").Value = 0.003236949171028675
$ws.Range("T2").Value = 0.003236949171028674
$ws.Range("G3").Value = 2.460980666666667
$ws.Range("H3").Value = 7.382942000000001
$ws.Range("I3").Value = 0.3244396275423151
$ws.Range("J3").Value = 0.3244396275423151
$ws.Range("O3").Value = 0.9569553279219795
$ws.Range("P3").Value = 0.9569553279219793
$ws.Range("Q3").Value = 104.9232962984936
$ws.Range("R3").Value = 944.3096666864421
$ws.Range("S3").Value = 0.3104742301656411
$ws.Range("T3").Value = 0.310474230165641
$ws.Range("G4").Value = 2.460980666666667
$ws.Range("H4").Value = 7.382942000000001
$ws.Range("I4").Value = 0.3244396275423151
$ws.Range("J4").Value = 0.3244396275423151
$ws.Range("M4").Value = 0.851471
$ws.Range("N4").Value = 2.554413
$ws.Range("O4").Value = 0.01911163320180161
$ws.Range("P4").Value = 0.01911163320180161
$ws.Range("Q4").Value = 2.095453669227334
$ws.Range("R4").Value = 18.859083023046
$ws.Range("S4").Value = 0.006200571157717859
$ws.Range("T4").Value = 0.006200571157717858
$ws.Range("G5").Value = 2.460980666666667
$ws.Range("H5").Value = 7.382942000000001
$ws.Range("I5").Value = 0.3244396275423151
$ws.Range("J5").Value = 0.3244396275423151
$ws.Range("M5").Value = 0.6217743333333333
$ws.Range("N5").Value = 1.865323
$ws.Range("O5").Value = 0.01395599262095996
$ws.Range("P5").Value = 0.01395599262095996
$ws.Range("Q5").Value = 1.530174613362889
$ws.Range("R5").Value = 13.771571520266
$ws.Range("S5").Value = 0.004527877047927548
$ws.Range("T5").Value = 0.004527877047927547
$ws.Range("I6").Value = 0.01481346816030475
$ws.Range("J6").Value = 0.01481346816030476
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4445023333333333
$ws.Range("N6").Value = 1.333507
$ws.Range("O6").Value = 0.009977046255258984
$ws.Range("P6").Value = 0.009977046255258982
$ws.Range("Q6").Value = 0.049946504685
$ws.Range("R6").Value = 0.449518542165
$ws.Range("S6").Value = 0.0001477946570361667
$ws.Range("T6").Value = 0.0001477946570361667
$ws.Range("I7").Value = 0.01481346816030475
$ws.Range("J7").Value = 0.01481346816030476
$ws.Range("O7").Value = 0.9569553279219795
$ws.Range("P7").Value = 0.9569553279219793
$ws.Range("S7").Value = 0.01417582728100624
$ws.Range("T7").Value = 0.01417582728100624
$ws.Range("I8").Value = 0.01481346816030475
$ws.Range("J8").Value = 0.01481346816030476
$ws.Range("M8").Value = 0.851471
$ws.Range("N8").Value = 2.554413
$ws.Range("O8").Value = 0.01911163320180161
$ws.Range("P8").Value = 0.01911163320180161
$ws.Range("Q8").Value = 0.095675538915
$ws.Range("R8").Value = 0.8610798502349999
$ws.Range("S8").Value = 0.0002831095699263114
$ws.Range("T8").Value = 0.0002831095699263114
$ws.Range("I9").Value = 0.01481346816030475
$ws.Range("J9").Value = 0.01481346816030476
$ws.Range("M9").Value = 0.6217743333333333
$ws.Range("N9").Value = 1.865323
$ws.Range("O9").Value = 0.01395599262095996
$ws.Range("P9").Value = 0.01395599262095996
$ws.Range("Q9").Value = 0.069865672965
$ws.Range("R9").Value = 0.628791056685
$ws.Range("S9").Value = 0.0002067366523360384
$ws.Range("T9").Value = 0.0002067366523360384
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.07607999999999999
$ws.Range("H10").Value = 0.22824
$ws.Range("I10").Value = 0.01002989060326601
$ws.Range("J10").Value = 0.01002989060326601
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4445023333333333
$ws.Range("N10").Value = 1.333507
$ws.Range("O10").Value = 0.009977046255258984
$ws.Range("P10").Value = 0.009977046255258982
$ws.Range("Q10").Value = 0.03381773751999999
$ws.Range("R10").Value = 0.30435963768
$ws.Range("S10").Value = 0.0001000686824839724
$ws.Range("T10").Value = 0.0001000686824839724
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.07607999999999999
$ws.Range("H11").Value = 0.22824
$ws.Range("I11").Value = 0.01002989060326601
$ws.Range("J11").Value = 0.01002989060326601
$ws.Range("O11").Value = 0.9569553279219795
$ws.Range("P11").Value = 0.9569553279219793
$ws.Range("Q11").Value = 3.24365180536
$ws.Range("R11").Value = 29.19286624824
$ws.Range("S11").Value = 0.009598157251270008
$ws.Range("T11").Value = 0.009598157251270008
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.07607999999999999
$ws.Range("H12").Value = 0.22824
$ws.Range("I12").Value = 0.01002989060326601
$ws.Range("J12").Value = 0.01002989060326601
$ws.Range("M12").Value = 0.851471
$ws.Range("N12").Value = 2.554413
$ws.Range("O12").Value = 0.01911163320180161
$ws.Range("P12").Value = 0.01911163320180161
$ws.Range("Q12").Value = 0.06477991368
$ws.Range("R12").Value = 0.58301922312
$ws.Range("S12").Value = 0.0001916875902638168
$ws.Range("T12").Value = 0.0001916875902638168
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.07607999999999999
$ws.Range("H13").Value = 0.22824
$ws.Range("I13").Value = 0.01002989060326601
$ws.Range("J13").Value = 0.01002989060326601
$ws.Range("M13").Value = 0.6217743333333333
$ws.Range("N13").Value = 1.865323
$ws.Range("O13").Value = 0.01395599262095996
$ws.Range("P13").Value = 0.01395599262095996
$ws.Range("Q13").Value = 0.04730459128
$ws.Range("R13").Value = 0.42574132152
$ws.Range("S13").Value = 0.0001399770792482161
$ws.Range("T13").Value = 0.0001399770792482161
$ws.Range("G14").Value = 4.935901333333334
$ws.Range("H14").Value = 14.807704
$ws.Range("I14").Value = 0.6507170136941141
$ws.Range("J14").Value = 0.6507170136941141
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4445023333333333
$ws.Range("N14").Value = 1.333507
$ws.Range("O14").Value = 0.009977046255258984
$ws.Range("P14").Value = 0.009977046255258982
$ws.Range("Q14").Value = 2.194019659769778
$ws.Range("R14").Value = 19.746176937928
$ws.Range("S14").Value = 0.006492233744710169
$ws.Range("T14").Value = 0.006492233744710169
$ws.Range("G15").Value = 4.935901333333334
$ws.Range("H15").Value = 14.807704
$ws.Range("I15").Value = 0.6507170136941141
$ws.Range("J15").Value = 0.6507170136941141
$ws.Range("O15").Value = 0.9569553279219795
$ws.Range("P15").Value = 0.9569553279219793
$ws.Range("Q15").Value = 210.4409210166338
$ws.Range("R15").Value = 1893.968289149704
$ws.Range("S15").Value = 0.6227071132240621
$ws.Range("T15").Value = 0.622707113224062
$ws.Range("G16").Value = 4.935901333333334
$ws.Range("H16").Value = 14.807704
$ws.Range("I16").Value = 0.6507170136941141
$ws.Range("J16").Value = 0.6507170136941141
$ws.Range("M16").Value = 0.851471
$ws.Range("N16").Value = 2.554413
$ws.Range("O16").Value = 0.01911163320180161
$ws.Range("P16").Value = 0.01911163320180161
$ws.Range("Q16").Value = 4.202776844194667
$ws.Range("R16").Value = 37.824991597752
$ws.Range("S16").Value = 0.01243626488389363
$ws.Range("T16").Value = 0.01243626488389362
$ws.Range("G17").Value = 4.935901333333334
$ws.Range("H17").Value = 14.807704
$ws.Range("I17").Value = 0.6507170136941141
$ws.Range("J17").Value = 0.6507170136941141
$ws.Range("M17").Value = 0.6217743333333333
$ws.Range("N17").Value = 1.865323
$ws.Range("O17").Value = 0.01395599262095996
$ws.Range("P17").Value = 0.01395599262095996
$ws.Range("Q17").Value = 3.069016760932445
$ws.Range("R17").Value = 27.621150848392
$ws.Range("S17").Value = 0.009081401841448157
$ws.Range("T17").Value = 0.009081401841448155
